# Page Object With Factories
# Applies the TestData.xlsx changes:
#  - TestSuite: add a FlightSearchTest/Y row
#  - LoginTest: rename RunMode header to lowercase runmode, add two new
#    login rows (ram@gmail.com / Hr@gail.com) with mailto hyperlinks
#  - CreateAccountTest sheet tab is renamed to FlightSearchTest (data is
#    unchanged)
#  - LoginTest becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

# --- Sheet 1: TestSuite -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A3").Value = "FlightSearchTest"
$ws1.Range("B3").Value = "Y"

# --- Sheet 2: LoginTest --------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# New login rows first (keeps shared-string append order matching the
# original authoring order), then fix up the RunMode -> runmode header.
$ws2.Range("A3").Value = "ram@gmail.com"
$ws2.Range("B3").Value = "123selenium*"
$ws2.Range("C3").Value = "Y"

$ws2.Range("C4").Value = "N"
$ws2.Range("A4").Value = "Hr@gail.com"
$ws2.Range("B4").Value = "123selenium*"

$ws2.Range("C1").Value = "runmode"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:ram@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:Hr@gail.com") | Out-Null

# --- Sheet 3: CreateAccountTest -> FlightSearchTest ----------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "FlightSearchTest"

# --- View state: selections + active tab ---------------------------------
# Select on each sheet in turn; whichever sheet is touched last ends up
# tabSelected / the workbook's active tab, so LoginTest must be last.
$ws1.Range("A3").Select() | Out-Null
$ws3.Range("D21").Select() | Out-Null
$ws2.Range("C1").Select() | Out-Null

Write-Host "TestData.xlsx updated: FlightSearchTest sheet + LoginTest rows"
